$d = $word.ActiveDocument

# Remove the inline picture ("A couple of men with blue hair...") from the
# first paragraph; this also removes the <w:r> run that wrapped the
# <w:drawing> element.
if ($d.InlineShapes.Count -gt 0) {
    $d.InlineShapes.Item(1).Delete()
}

# Tag the (now empty) paragraph mark's run properties with the
# English (Canada) language, matching the document's theme font language.
$p = $d.Paragraphs.Item(1)
$p.Range.LanguageID = "en-CA"
